$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B23").Value = 6321
$ws.Range("C23").Value = 1005
$ws.Range("D23").Value = 5916004
$ws.Range("E23").Value = 935.9284923271634
$ws.Range("F23").Value = 8.45916266300617
$ws.Range("G23").Value = 4.57856399583767
$ws.Range("H23").Value = 26.76482291399047
